$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so numeric-looking values are not
# auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '95.678.77'
$ws.Range('E2').Value = '  +2.06%  '

$ws.Range('D3').Value = '3.556.36'
$ws.Range('E3').Value = '  +7.00%  '

$ws.Range('E4').Value = '  -0.05%  '

$ws.Range('D5').Value = '238.66'
$ws.Range('E5').Value = '  +3.47%  '

$ws.Range('D6').Value = '635.89'
$ws.Range('E6').Value = '  +3.04%  '

$ws.Range('E7').Value = '  +6.82%  '

$ws.Range('D8').Value = '0.399'
$ws.Range('E8').Value = '  +3.39%  '

$ws.Range('D9').Value = '1.00'

$ws.Range('D10').Value = '1.01'
$ws.Range('E10').Value = '  +8.97%  '

$ws.Range('D11').Value = '3.552.49'
$ws.Range('E11').Value = '  +6.93%  '

$ws.Range('D12').Value = '42.84'
$ws.Range('E12').Value = '  +2.47%  '

$ws.Range('D13').Value = '0.200'
$ws.Range('E13').Value = '  +4.00%  '

$ws.Range('D14').Value = '6.40'
$ws.Range('E14').Value = '  +7.79%  '

$ws.Range('D15').Value = '4.231.06'
$ws.Range('E15').Value = '  +6.98%  '

$ws.Range('D16').Value = '95.629.73'
$ws.Range('E16').Value = '  +2.08%  '

$ws.Range('E17').Value = '  +4.13%  '

$ws.Range('D18').Value = '3.556.11'
$ws.Range('E18').Value = '  +6.97%  '

$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').Value = '12.94'
$ws.Range('E19').Value = '  +19.05%  '

$ws.Range('B20').Value = 'Polkadot'
$ws.Range('C20').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D20').Value = '7.93'
$ws.Range('E20').Value = '  -1.62%  '

$ws.Range('E21').Value = '  +3.99%  '

$ws.Range('D22').Value = '0.499'
$ws.Range('E22').Value = '  +10.91%  '

$ws.Range('D23').Value = '513.76'
$ws.Range('E23').Value = '  +4.49%  '

$ws.Range('E24').Value = '  -0.80%  '

$ws.Range('D25').Value = '6.69'
$ws.Range('E25').Value = '  +10.94%  '

$ws.Range('D26').Value = '0.0000192'
$ws.Range('E26').Value = '  +6.34%  '

$ws.Range('D27').Value = '96.13'
$ws.Range('E27').Value = '  +7.39%  '

$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').Value = '3.742.25'
$ws.Range('E28').Value = '  +6.61%  '

$ws.Range('B29').Value = 'Aptos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D29').Value = '12.29'
$ws.Range('E29').Value = '  +5.57%  '

$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').Value = '3.05'
$ws.Range('E30').Value = '  +16.40%  '

$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').Value = '0.145'
$ws.Range('E31').Value = '  +6.22%  '

$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = '11.47'
$ws.Range('E32').Value = '  +3.84%  '

$ws.Range('B33').Value = 'Dai'
$ws.Range('C33').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D33').Value = '0.999'
$ws.Range('E33').Value = '  -0.04%  '

$ws.Range('B34').Value = 'Cronos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D34').Value = '0.183'
$ws.Range('E34').Value = '  +5.71%  '

$ws.Range('B35').Value = 'Binance-PegBSC-USD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  +1.23%  '

$ws.Range('B36').Value = 'EthereumClassic'
$ws.Range('C36').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D36').Value = '30.00'
$ws.Range('E36').Value = '  +6.35%  '

$ws.Range('B37').Value = 'PolygonEcosystemToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D37').Value = '0.564'
$ws.Range('E37').Value = '  +6.79%  '

$ws.Range('B38').Value = 'Bittensor'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D38').Value = '578.34'
$ws.Range('E38').Value = '  +9.69%  '

$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D39').Value = '7.74'
$ws.Range('E39').Value = '  +5.47%  '

$ws.Range('B40').Value = 'Fetch.AI'
$ws.Range('C40').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D40').Value = '1.45'
$ws.Range('E40').Value = '  +7.20%  '

$ws.Range('B41').Value = 'USDe'
$ws.Range('C41').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  -0.05%  '

$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').Value = '0.151'
$ws.Range('E42').Value = '  +2.50%  '

$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').Value = '0.922'
$ws.Range('E43').Value = '  +6.57%  '

$ws.Range('B44').Value = 'ImmutableX'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D44').Value = '1.74'
$ws.Range('E44').Value = '  +3.97%  '

$ws.Range('B45').Value = 'WhiteBITCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D45').Value = '23.84'
$ws.Range('E45').Value = '  -0.96%  '

$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').Value = '0.0427'
$ws.Range('E46').Value = '  +3.46%  '

$ws.Range('B47').Value = 'Filecoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D47').Value = '5.61'
$ws.Range('E47').Value = '  +4.53%  '

$ws.Range('B48').Value = 'MantraDAO'
$ws.Range('C48').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D48').Value = '3.55'
$ws.Range('E48').Value = '  -5.13%  '

$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').Value = '2.17'
$ws.Range('E49').Value = '  +3.39%  '

$ws.Range('B50').Value = 'OKB'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D50').Value = '53.86'
$ws.Range('E50').Value = '  +1.80%  '

$ws.Range('B51').Value = 'Cosmos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D51').Value = '8.15'
$ws.Range('E51').Value = '  +2.59%  '

# Restore default (unstyled) formatting for column D now that values are set,
# matching the workbook's original style (no explicit number format).
$ws.Range("D2:D51").Style = "Normal"